# TC08_C3DC_phs002529_DiseasePhase-Unknown.xlsx
# "Updated remaining queries for C3DC"
#
# The stored DuckDB-style SQL queries (StatQuery in C2, plus the TabQuery
# column B2:B7) were joining on the raw primary-key columns ("id") instead
# of the renamed ("<table>_id") columns. Update every query's LEFT JOIN
# block so it matches the renamed schema, then restore the normal
# selection/scroll state and widen column C to fit the new text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldJoin = "LEFT JOIN `n    df_participant prt ON std.id = prt.`"study.id`"`nLEFT JOIN `n    df_diagnoses dgn ON prt.id = dgn.`"participant.id`"`nLEFT JOIN `n    df_treatments trt ON prt.id = trt.`"participant.id`"`nLEFT JOIN `n    df_treatment_resp trr ON prt.id = trr.`"participant.id`"`nLEFT JOIN `n    df_survival srv ON prt.id = srv.`"participant.id`"`nLEFT JOIN `n    df_reference_files rfs ON std.id = rfs.`"study.id`""
$newJoin = "LEFT JOIN `n    df_participant prt ON std.study_id = prt.`"study.study_id`"`nLEFT JOIN `n    df_diagnoses dgn ON prt.participant_id = dgn.`"participant.participant_id`"`nLEFT JOIN `n    df_treatments trt ON prt.participant_id = trt.`"participant.participant_id`"`nLEFT JOIN `n    df_treatment_resp trr ON prt.participant_id = trr.`"participant.participant_id`"`nLEFT JOIN `n    df_survival srv ON prt.participant_id = srv.`"participant.participant_id`"`nLEFT JOIN `n    df_reference_files rfs ON std.study_id = rfs.`"study.study_id`""

$queryCells = @("C2", "B2", "B3", "B4", "B5", "B6", "B7")
foreach ($cellRef in $queryCells) {
    $cell = $ws.Range($cellRef)
    $text = $cell.Value()
    if ($text.Contains($oldJoin)) {
        $cell.Value = $text.Replace($oldJoin, $newJoin)
    }
}

# Column C's text got a bit longer after the edit - stop relying on the old
# best-fit width and set an explicit fixed width (68 display characters).
$ws.Columns("C").ColumnWidth = 67.16666666666667

# Leave the sheet scrolled to the top with B2 selected instead of the old
# scrolled-down C7 selection.
$ws.Range("B2").Select()
